$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "Ось 3 поради, як подбати про безпеку дитини, навіть коли навколо небезпечно:" "Ось 3 поради, поради, як подбати про безпеку ваших дітей у небезпечній ситуації:"

Replace-Text "Пояснюйте дитині, що відбувається, простою мовою, зрозумілою для її віку." "Пояснюйте дітям, що відбувається, простими словами, зрозуміло для їхнього віку."

Replace-Text "Говоріть також про те, що може бути небезпечним." "Поговоріть з ними про речі, які можуть бути небезпечними."

Replace-Text "Заспокоюйте" "Підтримуйте"

Replace-Text "Ось дві речі, які допоможуть вам і вашим дітям бути в безпеці під час подорожі:" "Ось дві поради, які допоможуть вам і вашим дітям бути в безпеці під час подорожі:"

Replace-Text "Сьогоднішні поговоримо про те, як захистити дитину від людей, які можуть заманити її в небезпеку." "Сьогоднішні поговоримо про те, як захистити дитину від людей, які можуть завдати їй шкоду."

Replace-Text "Вчіться помічати ознаки, що щось може бути не так." "Навчіться помічати ознаки, які можуть свідчити про те, що щось не так."

Replace-Text "Зробіть укриття знайомим і затишним " "Зробіть укриття більш знайомим "

Replace-Text "Знаходьте разом щось, що викликає усмішку. " "Знайдіть щось, що змусить вас усміхнутись разом. Навіть маленька радість може зробити день світлішим. "

Replace-Text "Наскільки можливо, заспокоюйте дітей і говоріть їм, що все буде добре. " "Намагайтеся якомога частіше заспокоювати своїх дітей. Говоріть їм, що все буде добре. "

Replace-Text "Щодо маленьких дітей - обійміть їх або потримайте на руках, коли їм страшно. " "Для молодших дітей іноді найкраща підтримка — це обійми чи просто бути поруч, коли їм страшно. "

Replace-Text "Щодо старших - дайте знати, що ви поруч і готові вислухати, залишайте їм простір для розмови." "Для старших дітей важливо дати зрозуміти, що ви поруч і готові вислухати, коли вони будуть готові поділитися. Дайте їм простір, але не залишайте без підтримки."
